$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert new column before D, shifting the rest right
$ws.Columns("D:D").Insert()

# copy formatting from new E column (old D) into new D column, limited range
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)  # xlPasteFormats = -4122
$excel.CutCopyMode = $false

Write-Host $ws.Range("D7").NumberFormat
Write-Host $ws.Range("D8").NumberFormat
Write-Host $ws.Range("D7").Font.Name
Write-Host $ws.UsedRange.Address
